$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Borders.LineStyle = -4142
$ws.Range("S7").Value = 36
